$wb = $excel.ActiveWorkbook

$wsInputs = $wb.Worksheets.Item("Inputs")
$wsInputs.Range("C1").Value = Get-Date -Year 2019 -Month 7 -Day 1 -Hour 0 -Minute 0 -Second 0
$wsInputs.Range("C3").Value = "Monthly"
$wsInputs.Range("C4").Value = "No"

$wsPrices = $wb.Worksheets.Item("Prices")
$wsPrices.Range("C29").Formula = '=$L$28*IF(Inputs!$C$4 = "Yes", 6750, 9600)'

$excel.Calculate()
